# Apply "adding averages and more checks" update
# - refresh PERIOD TO EXPIRE (H) and LAST UPDATE (I) values across Training Dashboard
# - row 18 (LOTO SOPs) training becomes expired -> restyle as NOT VALID + new values
# - restyle title/header fonts (bold white) on both sheets
# - update Exam Dashboard comments + column width

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)   # Exam Dashboard

# ---------------------------------------------------------------------------
# 1. Training Dashboard: PERIOD TO EXPIRE (col H) + LAST UPDATE (col I)
# ---------------------------------------------------------------------------

$periodUpdates = @{
  3  = 386
  4  = 382
  5  = 360
  6  = 358
  7  = 405
  8  = 406
  9  = 342
  10 = 344
  11 = 348
  12 = 446
  13 = 409
  14 = 408
  15 = 386
  16 = 426
  17 = 427
  19 = -104
  20 = -190
  21 = -45
}

foreach ($row in $periodUpdates.Keys) {
    $ws1.Range("H$row").Value = $periodUpdates[$row]
    # use a text formula so the literal string is preserved instead of being
    # auto-converted to a date serial number by the engine
    $ws1.Range("I$row").Formula = '="16-Sep-2025"'
}

# ---------------------------------------------------------------------------
# 2. Row 18 (LOTO (SOPs)) training has now expired -> becomes NOT VALID
#    Copy the "NOT VALID" look (pink fill) from row 19 onto row 18, then
#    refresh its values.
# ---------------------------------------------------------------------------

$ws1.Range("A19:K19").Copy() | Out-Null
$ws1.Range("A18:K18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws1.Range("H18").Value = 14
$ws1.Range("I18").Formula = '="16-Sep-2025"'
$ws1.Range("J18").Value = "NOT VALID"

# ---------------------------------------------------------------------------
# 3. Title + header fonts: bold white text (instead of bold black/size-14)
# ---------------------------------------------------------------------------

$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215

$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

# ---------------------------------------------------------------------------
# 4. Exam Dashboard: comments + column width
# ---------------------------------------------------------------------------

$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"

$ws2.Columns.Item(5).ColumnWidth = $ws2.Columns.Item(2).ColumnWidth

Write-Output "edit complete"
